# "Scrape Current Makert Data"
# - Rename Sheet1 -> "List of Available Parameters"
# - Highlight the static header row (A1:H1) on MarketData with the
#   "Purple, Accent 4, Lighter 60%" theme fill
# - Add a note on H1 explaining the static columns
# - Restore the selections that were left on each sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Rename the second sheet
$ws2.Name = "List of Available Parameters"

# 2. Fill the static header columns A1:H1 with the Accent4/Lighter 60% color
#    (theme color 8 / Accent4, tint 0.6 -> RGB CC,C1,DA -> OLE BGR long)
$headerRange = $ws1.Range("A1:H1")
$headerRange.Interior.Color = 14336460

# 3. Comment on H1 describing the static columns
$commentText = "nobin thomas:" + [char]10 + "Static Colmns.Update the next colmn headers to fetch the required data"
$ws1.Range("H1").AddComment($commentText)

# 4. Leave the selections where the author left them: D22 on the parameters
#    sheet, then back to I8 on MarketData (which stays the active tab)
$ws2.Activate()
$ws2.Range("D22").Select()
$ws1.Activate()
$ws1.Range("I8").Select()
